$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($sheet, $addr, $val) {
    $r = $sheet.Range($addr)
    $r.NumberFormat = "@"
    $r.Value = $val
    $r.ClearFormats()
}

Set-TextValue $ws "D2" "41.496.44"
Set-TextValue $ws "E2" "  -0.80%  "
Set-TextValue $ws "D3" "2.464.67"
Set-TextValue $ws "E3" "  -0.56%  "
Set-TextValue $ws "D4" "1.00"
Set-TextValue $ws "E4" "  +0.36%  "
Set-TextValue $ws "D5" "311.16"
Set-TextValue $ws "E5" "  +0.13%  "
Set-TextValue $ws "D6" "90.86"
Set-TextValue $ws "E6" "  -4.09%  "
Set-TextValue $ws "D7" "0.535"
Set-TextValue $ws "E7" "  -3.33%  "
Set-TextValue $ws "E8" "  +0.33%  "
Set-TextValue $ws "E9" "  -4.42%  "
Set-TextValue $ws "D10" "31.86"
Set-TextValue $ws "E10" "  -6.20%  "
Set-TextValue $ws "E11" "  -1.69%  "
Set-TextValue $ws "E12" "  +0.59%  "
Set-TextValue $ws "D13" "2.840.81"
Set-TextValue $ws "E13" "  -1.06%  "
Set-TextValue $ws "E14" "  -3.87%  "
Set-TextValue $ws "D15" "15.15"
Set-TextValue $ws "E15" "  +2.82%  "
Set-TextValue $ws "D16" "2.417.68"
Set-TextValue $ws "E16" "  -0.91%  "
Set-TextValue $ws "D17" "0.759"
Set-TextValue $ws "E17" "  -3.75%  "
Set-TextValue $ws "D18" "41.291.70"
Set-TextValue $ws "E18" "  -1.33%  "
Set-TextValue $ws "D19" "6.19"
Set-TextValue $ws "E19" "  -3.17%  "
Set-TextValue $ws "D20" "0.0₃0908"
Set-TextValue $ws "E20" "  -1.19%  "
Set-TextValue $ws "D21" "69.86"
Set-TextValue $ws "E21" "  +0.68%  "
Set-TextValue $ws "E22" "  -7.09%  "
Set-TextValue $ws "D23" "232.30"
Set-TextValue $ws "E23" "  -1.80%  "
Set-TextValue $ws "E24" "  -4.41%  "
Set-TextValue $ws "E25" "  +0.19%  "
Set-TextValue $ws "E26" "  -4.07%  "
Set-TextValue $ws "D27" "23.78"
Set-TextValue $ws "E27" "  -3.90%  "
Set-TextValue $ws "E28" "  +0.32%  "
Set-TextValue $ws "E29" "  -2.26%  "
Set-TextValue $ws "D30" "35.48"
Set-TextValue $ws "E30" "  -2.77%  "
Set-TextValue $ws "D31" "151.27"
Set-TextValue $ws "E31" "  -2.19%  "
Set-TextValue $ws "D32" "5.31"
Set-TextValue $ws "E32" "  -5.28%  "
Set-TextValue $ws "D33" "2.54"
Set-TextValue $ws "E33" "  -3.34%  "
Set-TextValue $ws "D34" "0.0750"
Set-TextValue $ws "E34" "  -0.84%  "
Set-TextValue $ws "D35" "17.61"
Set-TextValue $ws "E35" "  +2.58%  "
Set-TextValue $ws "D36" "2.47"
Set-TextValue $ws "E36" "  -3.47%  "
Set-TextValue $ws "D37" "2.91"
Set-TextValue $ws "E37" "  -3.71%  "
Set-TextValue $ws "D38" "1.80"
Set-TextValue $ws "E38" "  -4.48%  "
Set-TextValue $ws "E39" "  -2.61%  "
Set-TextValue $ws "D40" "0.0994"
Set-TextValue $ws "E40" "  -7.10%  "
Set-TextValue $ws "D41" "4.02"
Set-TextValue $ws "E41" "  +0.10%  "
Set-TextValue $ws "E42" "  +0.70%  "
Set-TextValue $ws "D43" "20.19"
Set-TextValue $ws "E43" "  -3.81%  "
Set-TextValue $ws "D44" "1.936.46"
Set-TextValue $ws "E44" "  -2.99%  "
Set-TextValue $ws "E45" "  -3.28%  "
Set-TextValue $ws "D46" "2.89"
Set-TextValue $ws "E46" "  -5.76%  "
Set-TextValue $ws "D47" "8.62"
Set-TextValue $ws "E47" "  -0.57%  "
Set-TextValue $ws "D48" "2.692.05"
Set-TextValue $ws "E48" "  -1.24%  "
Set-TextValue $ws "D49" "94.30"
Set-TextValue $ws "E49" "  -3.48%  "
Set-TextValue $ws "E50" "  -4.72%  "
Set-TextValue $ws "D51" "65.67"
Set-TextValue $ws "E51" "  -5.01%  "
